$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns affected by this update (indices: B=2,C=3,E=5,F=6,G=7,H=8,J=10,L=12,N=14)
$cols = 2, 3, 5, 6, 7, 8, 10, 12, 14

$data = @(
    @(1.829825935648728, 0.1912749938038871, 0.1748610333157306, 2.291731066012929, 1.229822858521331, 1.152188506180579, 0.0575677407483699, 0.4519815462856656, 1.482805291469848),
    @(1.717763035893597, 0.175700259572011, 0.1752132567747999, 2.281435561401068, 1.218442010815423, 1.153133832300298, 0.05827408791961242, 0.4434309366823044, 1.502454042642874),
    @(1.649730772083615, 0.166058391701057, 0.17546805021572, 2.276419841258814, 1.21237134554886, 1.154368033702241, 0.0587291881566685, 0.4383570897230555, 1.515126758487728),
    @(1.622202356267564, 0.1621093399897973, 0.1755815969999173, 2.274703882269108, 1.210127316077433, 1.155035131039469, 0.05892002478969705, 0.4363339330049456, 1.520443832529398),
    @(1.617643092203537, 0.161452398375701, 0.1756010390435705, 2.274438747946235, 1.209768552916074, 1.155155809217106, 0.05895203793649628, 0.4360006802647263, 1.521335956231068),
    @(1.649358722938302, 0.1660052141194797, 0.1754695421676953, 2.276395371704169, 1.212340152365456, 1.154376366103293, 0.05873174006712967, 0.4383296244307644, 1.515197847735203),
    @(1.79102606477943, 0.1859211595847228, 0.1749745006734944, 2.287909939721359, 1.225707867152295, 1.152378614155097, 0.05780684667039715, 0.4489968017343671, 1.489453854471625),
    @(2.074979770836649, 0.2243549131588622, 0.1743082175549269, 2.320871219521564, 1.259240253437156, 1.153660909290977, 0.05616298511558515, 0.4713084011461746, 1.443801406876764),
    @(2.287366520623607, 0.2522226299761314, 0.1740027862676055, 2.351451651404261, 1.288398922908613, 1.157792494192165, 0.05505882075855162, 0.4885453162095956, 1.413211571197863),
    @(2.384810558324432, 0.264822283381136, 0.1739034920723661, 2.366753389008593, 1.3026594611799, 1.160369025550921, 0.05457900963863871, 0.4965693752637037, 1.399937525927065),
    @(2.421829168160514, 0.2695824164712803, 0.1738715662454364, 2.372748260148654, 1.308203807034261, 1.161445254622492, 0.0544005528584206, 0.4996340617421708, 1.39500338401135),
    @(2.413851277129879, 0.2685577284936471, 0.1738781900800479, 2.3714482353205, 1.307003304896853, 1.161208992236936, 0.05443884266000865, 0.4989728660644204, 1.396061924483478),
    @(2.387853725781667, 0.2652141241287893, 0.1739007519371025, 2.367242570894192, 1.303112702616318, 1.160455550510449, 0.05456426298311357, 0.496820985651766, 1.399529736994435),
    @(2.371944916368761, 0.2631646285983322, 0.1739153099570192, 2.364692596823772, 1.300748402400018, 1.160007150091644, 0.05464150828988235, 0.4955062977315805, 1.401665918616049),
    @(2.281014952591875, 0.2513976586940601, 0.1740100709173209, 2.350479668781176, 1.287487088576597, 1.157638165318502, 0.05509063028358963, 0.4880245925228479, 1.414091986596548),
    @(2.22544418470568, 0.2441591963126655, 0.1740783395116612, 2.342116986862251, 1.279607504301168, 1.156363626258127, 0.05537191280982512, 0.4834815456094361, 1.421879411856846),
    @(2.193559331570384, 0.239988524612869, 0.1741213392405374, 2.337437849318718, 1.275169077605057, 1.155696135806608, 0.05553581458323631, 0.4808857295443119, 1.426418893566815),
    @(2.182777075190245, 0.238575151211819, 0.174136540245577, 2.335876033947386, 1.273682367113963, 1.155481390665074, 0.05559167210292149, 0.4800097937629459, 1.427966244625548),
    @(2.231351720083296, 0.2449304980052602, 0.1740706859996219, 2.342993662137758, 1.280436594882417, 1.156492512464553, 0.05534175082514547, 0.4839633791528826, 1.421044178530176),
    @(2.395486619567748, 0.2661965220340221, 0.1738939711733387, 2.368472432900234, 1.304251545954429, 1.160674123093088, 0.05452733609497162, 0.4974523366008725, 1.398508645562146),
    @(2.503449760463809, 0.2800305331720381, 0.173811542282678, 2.386292797888828, 1.320656885137993, 1.163993294735178, 0.05401394075752819, 0.5064205000034718, 1.384319351140375),
    @(2.4457646672397, 0.2726529473005996, 0.1738525194705005, 2.376674655408891, 1.311823789807249, 1.162168042847583, 0.0542862210738706, 0.5016201319799762, 1.39184305534692),
    @(2.228680725445201, 0.244581821122182, 0.1740741344664443, 2.342596916458618, 1.280061477619284, 1.156434039756078, 0.05535538024087217, 0.4837454922293887, 1.421421593199831),
    @(1.997502937450463, 0.2140231357799678, 0.1744560267406037, 2.310839331686367, 1.249379095770308, 1.152755330557738, 0.05658950679088681, 0.4651238630860917, 1.455633805692926)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($r, $cols[$j]).Value = $rowValues[$j]
    }
}
